# Apply the "Add files via upload" edit to the purchase-order template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (Hoja1 -> HOJA1) ---
$ws.Name = "HOJA1"

# --- Remove the SUCURSAL dropdown data validation on B2 ---
$ws.Range("B2").Validation.Delete()

# --- Remove the logo picture that used to sit in the merged D1:F4 box ---
foreach ($shp in @($ws.Shapes)) {
    $shp.Delete()
}

# --- Re-format the now-empty logo placeholder cells.
#     D1 keeps the bold-ish Arial look with full centering (same format as
#     the B3/B4 "Arial 12, centered" style already used elsewhere).
$ws.Range("B3").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

#     E1:F4 (minus D1) become plain default-font cells that are only
#     vertically centered.
$ws.Range("E1").Style = "Normal"
$ws.Range("E1").VerticalAlignment = -4108   # xlCenter
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("D2:F4").PasteSpecial(-4122)

# --- Normalize C3 / C4 back to the plain style used by A1 (drop the stray
#     empty-alignment formatting they used to carry). ---
$ws.Range("A1").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Fill in the first order line (row 10) ---
$ws.Range("A10").Value = "V01009"
$ws.Range("B10").Value = "CILONAXOL 5/30 MG CAJA C/10 TABS"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "N/A"

# --- Reset the view: scroll back to A1, select D5 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D5").Select()
